$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the balance label in A9
$ws.Range("A9").Value = "BALANCE : 113.0"

# Append new transaction rows (20-28)
$data = @(
    @("2025-10-21 18:54:38", "Deposit",  10000,    109874),
    @("2025-10-21 19:00:17", "Deposit",  234,      110108),
    @("2025-10-21 19:01:33", "Deposit",  3,        110111),
    @("2025-10-21 19:06:41", "Deposit",  1,        110112),
    @("2025-10-21 19:06:47", "Deposit",  1,        110113),
    @("2025-10-21 19:36:06", "Withdraw", 100000,   10113),
    @("2025-10-21 19:36:34", "Withdraw", 10000,    113),
    @("2025-10-21 19:54:59", "Deposit",  10000000, 10000113),
    @("2025-10-21 19:57:16", "Withdraw", 10000000, 113)
)

$row = 20
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    # The source rows for this trailing transaction log carry no explicit
    # cell styling (unlike the column defaults), so strip the inherited
    # column format back off each new row.
    $ws.Rows($row).ClearFormats()
    $row++
}
